$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns C:H (epoch100 .. epoch350) from row 1 and row 2
$ws.Range("C1:H2").ClearContents()

# Update remaining value in B2
$ws.Range("B2").Value = 84.68468487262726
